$d = $word.ActiveDocument

# 1. Replace the "{{title}}" placeholder (originally 4 separate runs: "{{",
#    "title", "}", "}") in the Title paragraph with the literal text
#    "Eksport fra MIME" (plus a one-character placeholder "X" used below).
#    Find/Replace merges the matched runs into a single new run and keeps
#    the original run formatting (rFonts Arial, lang en-US).
$d.Content.Find.Execute("{{title}}", $false, $false, $false, $false, $false, $true, 1, $false, "Eksport fra MIMEX", 2) | Out-Null

# 2. Move the "_GoBack" bookmark so that it sits right after the new title
#    text instead of its old location further down (just before
#    "Hjemmel:"). We wrap a bookmark named "_GoBack" around the trailing
#    placeholder character "X" -- adding a bookmark whose name already
#    exists relocates the existing bookmark to the new spot -- and then
#    delete that placeholder character. Once its single character is gone,
#    the bookmark collapses to a zero-length bookmark sitting exactly where
#    the placeholder used to be, i.e. right after "Eksport fra MIME" and
#    before the paragraph mark.
$titlePara = $d.Paragraphs(1).Range
$placeholder = $d.Range($titlePara.End - 2, $titlePara.End - 1)
$d.Bookmarks.Add("_GoBack", $placeholder)
$placeholder.Text = ""
